$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 5051754
$ws.Range("I33").Value = 1519.8572
$ws.Range("J33").Value = 22727572
$ws.Range("K33").Value = 1519.8572
$ws.Range("L33").Value = 22727572
$ws.Range("M33").Value = -1290.8572
$ws.Range("N33").Value = -22728030
$ws.Range("H75").Value = 40000
$ws.Range("J75").Value = 40000
$ws.Range("L75").Value = 40000
$ws.Range("N75").Value = -41872
$ws.Range("H78").Value = 40000
$ws.Range("J78").Value = 40000
$ws.Range("L78").Value = 120000
$ws.Range("N78").Value = -129360
$ws.Range("H137").Value = 1574.4642
$ws.Range("I137").Value = 1204.5652
$ws.Range("K137").Value = 3613.6956
$ws.Range("M137").Value = -1063.6956

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4176.8237
$ws.Range("I32").Value = 2950.2036
$ws.Range("K32").Value = 2950.2036
$ws.Range("M32").Value = -2663.2036
$ws.Range("H76").Value = 29500
$ws.Range("J76").Value = 29500
$ws.Range("L76").Value = 29500
$ws.Range("N76").Value = -30176
$ws.Range("H79").Value = 29500
$ws.Range("J79").Value = 29500
$ws.Range("L79").Value = 29500
$ws.Range("N79").Value = -31840

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1646325.5
$ws.Range("I7").Value = 2006119.2
$ws.Range("J7").Value = 1046669.3
$ws.Range("K7").Value = 2006119.2
$ws.Range("L7").Value = 1046669.3
$ws.Range("M7").Value = -2006006.2
$ws.Range("N7").Value = -1046895.3
$ws.Range("H88").Value = 34000
$ws.Range("J88").Value = 34000
$ws.Range("L88").Value = 34000
$ws.Range("N88").Value = -34812
$ws.Range("H91").Value = 34000
$ws.Range("J91").Value = 34000
$ws.Range("L91").Value = 34000
$ws.Range("N91").Value = -36808

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 972
$ws.Range("I16").Value = 680.4286
$ws.Range("J16").Value = 3013
$ws.Range("K16").Value = 680.4286
$ws.Range("L16").Value = 3013
$ws.Range("M16").Value = -393.4286
$ws.Range("N16").Value = -3587
$ws.Range("H23").Value = 12670
$ws.Range("I23").Value = 12670
$ws.Range("K23").Value = 12670
$ws.Range("M23").Value = -12430
$ws.Range("H27").Value = 12670
$ws.Range("I27").Value = 12670
$ws.Range("K27").Value = 12670
$ws.Range("M27").Value = -12478
$ws.Range("H88").Value = 34227.668
$ws.Range("J88").Value = 37011
$ws.Range("L88").Value = 37011
$ws.Range("N88").Value = -37823
$ws.Range("H91").Value = 34227.668
$ws.Range("J91").Value = 37011
$ws.Range("L91").Value = 37011
$ws.Range("N91").Value = -39819
$ws.Range("H113").Value = 972
$ws.Range("I113").Value = 680.4286
$ws.Range("J113").Value = 3013
$ws.Range("K113").Value = 680.4286
$ws.Range("L113").Value = 3013
$ws.Range("M113").Value = 1489.5714
$ws.Range("N113").Value = -7353
$ws.Range("H134").Value = 2050.647
$ws.Range("I134").Value = 2105.5833
$ws.Range("J134").Value = 1918.8
$ws.Range("K134").Value = 6316.749899999999
$ws.Range("L134").Value = 5756.4
$ws.Range("M134").Value = -3781.749899999999
$ws.Range("N134").Value = -10826.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 353896.94
$ws.Range("I5").Value = 599.75
$ws.Range("J5").Value = 667938.9
$ws.Range("K5").Value = 1799.25
$ws.Range("L5").Value = 2003816.7
$ws.Range("M5").Value = -1687.25
$ws.Range("N5").Value = -2004040.7
$ws.Range("H22").Value = 2981.818
$ws.Range("I22").Value = 3100
$ws.Range("J22").Value = 2970
$ws.Range("K22").Value = 9300
$ws.Range("L22").Value = 8910
$ws.Range("M22").Value = -9131
$ws.Range("N22").Value = -9248
$ws.Range("H27").Value = 2981.818
$ws.Range("I27").Value = 3100
$ws.Range("J27").Value = 2970
$ws.Range("K27").Value = 9300
$ws.Range("L27").Value = 8910
$ws.Range("M27").Value = -9198
$ws.Range("N27").Value = -9114
$ws.Range("H39").Value = 3259.8
$ws.Range("J39").Value = 3259.8
$ws.Range("L39").Value = 9779.400000000001
$ws.Range("N39").Value = -10367.4
$ws.Range("H49").Value = 662.9091
$ws.Range("J49").Value = 699.1111
$ws.Range("L49").Value = 2097.3333
$ws.Range("N49").Value = -2409.3333
$ws.Range("H58").Value = 3566.389
$ws.Range("I58").Value = 450.5
$ws.Range("J58").Value = 3955.875
$ws.Range("K58").Value = 1351.5
$ws.Range("L58").Value = 11867.625
$ws.Range("M58").Value = -1223.5
$ws.Range("N58").Value = -12123.625
$ws.Range("H132").Value = 1612161.1
$ws.Range("I132").Value = 1905.6364
$ws.Range("J132").Value = 1917554.4
$ws.Range("K132").Value = 17150.7276
$ws.Range("L132").Value = 17257989.6
$ws.Range("M132").Value = -14620.7276
$ws.Range("N132").Value = -17263049.6
$ws.Range("H135").Value = 353896.94
$ws.Range("I135").Value = 599.75
$ws.Range("J135").Value = 667938.9
$ws.Range("K135").Value = 5397.75
$ws.Range("L135").Value = 6011450.100000001
$ws.Range("M135").Value = -2862.75
$ws.Range("N135").Value = -6016520.100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 69753
$ws.Range("J4").Value = 69753
$ws.Range("L4").Value = 69753
$ws.Range("N4").Value = -69977

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5625
$ws.Range("I2").Value = 3000
$ws.Range("J2").Value = 6000
$ws.Range("K2").Value = 3000
$ws.Range("L2").Value = 6000
$ws.Range("M2").Value = -2888
$ws.Range("N2").Value = -6224
$ws.Range("H61").Value = 2462.5
$ws.Range("I61").Value = 2528.5715
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 2528.5715
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -2326.5715
$ws.Range("N61").Value = -2404
$ws.Range("H113").Value = 2462.5
$ws.Range("I113").Value = 2528.5715
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 2528.5715
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -358.5715
$ws.Range("N113").Value = -6340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 34650
$ws.Range("J80").Value = 34650
$ws.Range("L80").Value = 34650
$ws.Range("N80").Value = -36646
$ws.Range("H83").Value = 34650
$ws.Range("J83").Value = 34650
$ws.Range("L83").Value = 103950
$ws.Range("N83").Value = -113934

Write-Output "Applied 166 cell updates across 8 sheets"